# Updated with Apr 2 data
#
# - row 28 (Apr 1, 2020): fill in the previously-missing "tested" value (col B)
# - row 29 (Apr 2, 2020): brand-new day of data across the whole column range
# - refresh the active selection to match the author's final cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 28: backfill the "tested" figure that was missing before ---
$ws.Range("B28").Value = 2607

# --- row 29: new day of data (2020-04-02, serial 43923) ---
$ws.Range("A29").Value = 43923
$ws.Range("C29").Value = 966
$ws.Range("D29").Value = 8
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 9
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 187
$ws.Range("I29").Value = 10
$ws.Range("J29").Value = 216
$ws.Range("K29").Value = 23
$ws.Range("L29").Value = 170
$ws.Range("M29").Value = 24
$ws.Range("N29").Value = 149
$ws.Range("O29").Value = 34
$ws.Range("P29").Value = 108
$ws.Range("Q29").Value = 32
$ws.Range("R29").Value = 71
$ws.Range("S29").Value = 27
$ws.Range("T29").Value = 45
$ws.Range("U29").Value = 29
$ws.Range("V29").Value = 3
$ws.Range("W29").Value = 0
$ws.Range("X29").Value = 439
$ws.Range("Y29").Value = 514
$ws.Range("Z29").Value = 13
$ws.Range("AA29").Value = 181
$ws.Range("AB29").Value = 70
$ws.Range("AC29").Value = 16
$ws.Range("AD29").Value = 32
$ws.Range("AE29").Value = 65
$ws.Range("AF29").Value = 2
$ws.Range("AG29").Value = 5
$ws.Range("AH29").Value = 43
$ws.Range("AI29").Value = 25
$ws.Range("AJ29").Value = 23
$ws.Range("AK29").Value = 1
$ws.Range("AL29").Value = 12
$ws.Range("AM29").Value = 8
$ws.Range("AN29").Value = 12
$ws.Range("AO29").Value = 22
$ws.Range("AP29").Value = 10
$ws.Range("AQ29").Value = 540
$ws.Range("AR29").Value = 12
$ws.Range("AS29").Value = 10
$ws.Range("AT29").Value = 5
$ws.Range("AU29").Value = 9
$ws.Range("AV29").Value = 6
$ws.Range("AW29").Value = 6
$ws.Range("AY29").Value = 8
$ws.Range("AZ29").Value = 2
$ws.Range("BA29").Value = 4
$ws.Range("BB29").Value = 14
$ws.Range("BD29").Value = 21
$ws.Range("BE29").Value = 40

# --- final view state: scroll so column O is leftmost, select S18 ---
$ws.Range("O4").Select()
$ws.Range("S18").Select()
